$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at the top of the data block (rows 2-7), pushing existing
# data (previously rows 2-21) down to rows 8-27.
$ws.Range("A2:A7").EntireRow.Insert()

# Row insertion inherits the header row's formatting by default; clear it so
# the new data rows look like the other (unstyled) data rows.
$ws.Range("A2:H7").ClearFormats()

# Populate the newly inserted rows 2-7 with the new sensor readings.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -3.130717563629151
$ws.Cells.Item(2, 4).Value = 5.486354422569275
$ws.Cells.Item(2, 5).Value = -2.054217553138733
$ws.Cells.Item(2, 6).Value = 0.0345138870179653
$ws.Cells.Item(2, 7).Value = -0.0546724386513233
$ws.Cells.Item(2, 8).Value = 0.1018617823719978

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -3.261763083934784
$ws.Cells.Item(3, 4).Value = 5.438373637199402
$ws.Cells.Item(3, 5).Value = -2.230073320865632
$ws.Cells.Item(3, 6).Value = -0.0178678091615438
$ws.Cells.Item(3, 7).Value = -0.0160352122038602
$ws.Cells.Item(3, 8).Value = 0.0375682115554809

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -3.289464282989502
$ws.Cells.Item(4, 4).Value = 5.444673538208008
$ws.Cells.Item(4, 5).Value = -2.207874870300293
$ws.Cells.Item(4, 6).Value = 0.0255036242306232
$ws.Cells.Item(4, 7).Value = -0.0097738439217209
$ws.Cells.Item(4, 8).Value = -0.0500909499824047

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -3.29152911901474
$ws.Cells.Item(5, 4).Value = 5.474263513088226
$ws.Cells.Item(5, 5).Value = -1.988579791784286
$ws.Cells.Item(5, 6).Value = 0.040775254368782
$ws.Cells.Item(5, 7).Value = 0.0166460778564214
$ws.Cells.Item(5, 8).Value = 0.0087048299610614

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -3.39021909236908
$ws.Cells.Item(6, 4).Value = 5.475549221038817
$ws.Cells.Item(6, 5).Value = -1.85912013053894
$ws.Cells.Item(6, 6).Value = 0.01328631862998
$ws.Cells.Item(6, 7).Value = 0.058643065392971
$ws.Cells.Item(6, 8).Value = 0.0251981914043426

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -3.51887332201004
$ws.Cells.Item(7, 4).Value = 5.510936594009399
$ws.Cells.Item(7, 5).Value = -1.821529471874237
$ws.Cells.Item(7, 6).Value = 0.0279470849782228
$ws.Cells.Item(7, 7).Value = 0.0377209298312664
$ws.Cells.Item(7, 8).Value = 0.0471893399953842

# Append 4 new rows (28-31) at the bottom with additional sensor readings.
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "falling"
$ws.Cells.Item(28, 3).Value = 1.438675880432129
$ws.Cells.Item(28, 4).Value = 5.703988456726075
$ws.Cells.Item(28, 5).Value = 1.146768474578857
$ws.Cells.Item(28, 6).Value = -0.0830776765942573
$ws.Cells.Item(28, 7).Value = -0.1788308024406433
$ws.Cells.Item(28, 8).Value = 0.0158824957907199

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = 1.517132639884949
$ws.Cells.Item(29, 4).Value = 5.866671967506409
$ws.Cells.Item(29, 5).Value = 1.15747617483139
$ws.Cells.Item(29, 6).Value = -0.0752891451120376
$ws.Cells.Item(29, 7).Value = 0.0606283769011497
$ws.Cells.Item(29, 8).Value = -0.0500909499824047

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = 1.545208883285523
$ws.Cells.Item(30, 4).Value = 5.865323352813721
$ws.Cells.Item(30, 5).Value = 1.205629134178161
$ws.Cells.Item(30, 6).Value = 0.0174096599221229
$ws.Cells.Item(30, 7).Value = -0.0296269636601209
$ws.Cells.Item(30, 8).Value = -0.0302378293126821

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = 1.573752522468567
$ws.Cells.Item(31, 4).Value = 5.771291553974152
$ws.Cells.Item(31, 5).Value = 1.240318953990936
$ws.Cells.Item(31, 6).Value = 0.0155770638957619
$ws.Cells.Item(31, 7).Value = 0.0525344125926494
$ws.Cells.Item(31, 8).Value = -0.0565050356090068
